$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 454, shifting rows 454:494 down to 455:495
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new record's data.
# Columns A,B,C,E,F,G,H,N,O,Q,R follow the same pattern as the surrounding rows.
$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 45106
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100112045
$ws.Cells.Item(454, 7).Value = "Zapallo"
$ws.Cells.Item(454, 8).Value = "Paine"
$ws.Cells.Item(454, 9).Value = "1a (guarda)"
$ws.Cells.Item(454, 10).Value = 500
$ws.Cells.Item(454, 11).Value = 550
$ws.Cells.Item(454, 12).Value = 600
$ws.Cells.Item(454, 13).Value = 575
$ws.Cells.Item(454, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(454, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(454, 16).Value = 575
$ws.Cells.Item(454, 17).Value = 1
$ws.Cells.Item(454, 18).Value = "Hortaliza"
